$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.185.46'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").Value = '3.744.23'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("D7").Value = '3.741.38'
$ws.Range("E7").Value = '  -1.38%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("E10").Value = '  -2.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.54'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("E12").Value = '  -2.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").Value = '4.369.07'
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").Value = '3.745.79'
$ws.Range("E16").Value = '  -1.46%  '
$ws.Range("D17").Value = '69.252.68'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("E18").Value = '  -2.81%  '
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '499.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.724'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("E26").Value = '  -3.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.94%  '
$ws.Range("E28").Value = '  -3.37%  '
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("E31").Value = '  +3.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.18%  '
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.10'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.347'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("E39").Value = '  +2.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '446.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.36%  '
$ws.Range("E41").Value = '  +9.58%  '
$ws.Range("E42").Value = '  -4.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.72'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("E45").Value = '  -2.64%  '
$ws.Range("D46").Value = '2.940.98'
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0359'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.24'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.79%  '
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.96%  '

Write-Output "edit complete"
